$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 685
$ws.Range("I2").Value = 1830
$ws.Range("J2").Value = 7792
$ws.Range("K2").Value = 43
$ws.Range("L2").Value = 2145
$ws.Range("M2").Value = 118
$ws.Range("N2").Value = 1356
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 19
$ws.Range("Q2").Value = 15
$ws.Range("R2").Value = 114
$ws.Range("S2").Value = 788
$ws.Range("T2").Value = 1468
$ws.Range("U2").Value = 117
$ws.Range("V2").Value = 11956
$ws.Range("W2").Value = 4
$ws.Range("X2").Value = 11738
$ws.Range("Y2").Value = 15
$ws.Range("Z2").Value = 188
$ws.Range("AA2").Value = 62
